# Scheduled runner update: refresh market-price-derived columns (H-N)
# on the Leve profit sheets. Values are sourced externally (Universalis)
# and written verbatim; there are no formulas on these sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3561.6775
$ws.Range("J80").Value = 7994.846
$ws.Range("L80").Value = 23984.538
$ws.Range("N80").Value = -25980.538
$ws.Range("H83").Value = 3561.6775
$ws.Range("J83").Value = 7994.846
$ws.Range("L83").Value = 71953.614
$ws.Range("N83").Value = -81937.614
$ws.Range("H115").Value = 550
$ws.Range("I115").Value = 308.33334
$ws.Range("K115").Value = 925.0000200000001
$ws.Range("M115").Value = 641.9999799999999
$ws.Range("H128").Value = 76778.336
$ws.Range("J128").Value = 76778.336
$ws.Range("L128").Value = 76778.336
$ws.Range("N128").Value = -86738.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6564.57
$ws.Range("I32").Value = 4663.047
$ws.Range("J32").Value = 17339.867
$ws.Range("K32").Value = 4663.047
$ws.Range("L32").Value = 17339.867
$ws.Range("M32").Value = -4376.047
$ws.Range("N32").Value = -17913.867
$ws.Range("H61").Value = 213100.42
$ws.Range("I61").Value = 5181.943
$ws.Range("J61").Value = 772880.9399999999
$ws.Range("K61").Value = 5181.943
$ws.Range("L61").Value = 772880.9399999999
$ws.Range("M61").Value = -4969.943
$ws.Range("N61").Value = -773304.9399999999
$ws.Range("H102").Value = 2648431.2
$ws.Range("I102").Value = 4116893
$ws.Range("J102").Value = 5200
$ws.Range("K102").Value = 4116893
$ws.Range("L102").Value = 5200
$ws.Range("M102").Value = -4115271
$ws.Range("N102").Value = -8444
$ws.Range("H133").Value = 25822.857
$ws.Range("J133").Value = 25822.857
$ws.Range("L133").Value = 25822.857
$ws.Range("N133").Value = -30882.857
$ws.Range("H136").Value = 213100.42
$ws.Range("I136").Value = 5181.943
$ws.Range("J136").Value = 772880.9399999999
$ws.Range("K136").Value = 15545.829
$ws.Range("L136").Value = 2318642.82
$ws.Range("M136").Value = -12995.829
$ws.Range("N136").Value = -2323742.82
$ws.Range("H138").Value = 42943
$ws.Range("J138").Value = 42943
$ws.Range("L138").Value = 42943
$ws.Range("N138").Value = -53223
$ws.Range("H139").Value = 59241.25
$ws.Range("J139").Value = 59241.25
$ws.Range("L139").Value = 59241.25
$ws.Range("N139").Value = -69521.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2114.2856
$ws.Range("I86").Value = 1900
$ws.Range("J86").Value = 2200
$ws.Range("K86").Value = 1900
$ws.Range("L86").Value = 2200
$ws.Range("M86").Value = -777
$ws.Range("N86").Value = -4446
$ws.Range("H89").Value = 2114.2856
$ws.Range("I89").Value = 1900
$ws.Range("J89").Value = 2200
$ws.Range("K89").Value = 9500
$ws.Range("L89").Value = 11000
$ws.Range("M89").Value = -3884
$ws.Range("N89").Value = -22232
$ws.Range("H94").Value = 1203.129
$ws.Range("I94").Value = 890.3333
$ws.Range("J94").Value = 1860
$ws.Range("K94").Value = 890.3333
$ws.Range("L94").Value = 1860
$ws.Range("M94").Value = -439.3333
$ws.Range("N94").Value = -2762
$ws.Range("H135").Value = 55192.727
$ws.Range("J135").Value = 55192.727
$ws.Range("L135").Value = 55192.727
$ws.Range("N135").Value = -65332.727
$ws.Range("H138").Value = 45312
$ws.Range("J138").Value = 45312
$ws.Range("L138").Value = 45312
$ws.Range("N138").Value = -55592

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1954.5385
$ws.Range("I16").Value = 1929.9412
$ws.Range("J16").Value = 2001
$ws.Range("K16").Value = 1929.9412
$ws.Range("L16").Value = 2001
$ws.Range("M16").Value = -1642.9412
$ws.Range("N16").Value = -2575
$ws.Range("H31").Value = 16114893
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 16114893
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 16114893
$ws.Range("M31").Value = ""
$ws.Range("N31").Value = -16115483
$ws.Range("H34").Value = 16114893
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 16114893
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 16114893
$ws.Range("M34").Value = ""
$ws.Range("N34").Value = -16115297
$ws.Range("H99").Value = 1152
$ws.Range("I99").Value = 1152
$ws.Range("K99").Value = 1152
$ws.Range("M99").Value = 346
$ws.Range("H113").Value = 1954.5385
$ws.Range("I113").Value = 1929.9412
$ws.Range("J113").Value = 2001
$ws.Range("K113").Value = 1929.9412
$ws.Range("L113").Value = 2001
$ws.Range("M113").Value = 240.0588
$ws.Range("N113").Value = -6341
$ws.Range("H126").Value = 1152
$ws.Range("I126").Value = 1152
$ws.Range("K126").Value = 3456
$ws.Range("M126").Value = -986

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 826.95386
$ws.Range("I92").Value = 768.44446
$ws.Range("K92").Value = 2305.33338
$ws.Range("M92").Value = -1057.33338
$ws.Range("H125").Value = 5610.8887
$ws.Range("J125").Value = 6124.75
$ws.Range("L125").Value = 18374.25
$ws.Range("N125").Value = -28214.25
$ws.Range("H131").Value = 1924254.4
$ws.Range("I131").Value = 8333774
$ws.Range("J131").Value = 1398.5
$ws.Range("K131").Value = 25001322
$ws.Range("L131").Value = 4195.5
$ws.Range("M131").Value = -24996282
$ws.Range("N131").Value = -14275.5
$ws.Range("H132").Value = 8716983
$ws.Range("I132").Value = 2268
$ws.Range("J132").Value = 10895662
$ws.Range("K132").Value = 20412
$ws.Range("L132").Value = 98060958
$ws.Range("M132").Value = -17882
$ws.Range("N132").Value = -98066018

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 45217.332
$ws.Range("J51").Value = 50326
$ws.Range("L51").Value = 50326
$ws.Range("N51").Value = -51344
$ws.Range("H102").Value = 1928.7273
$ws.Range("I102").Value = 1592.6666
$ws.Range("J102").Value = 2648.8572
$ws.Range("K102").Value = 1592.6666
$ws.Range("L102").Value = 2648.8572
$ws.Range("M102").Value = 29.33339999999998
$ws.Range("N102").Value = -5892.8572
$ws.Range("H141").Value = 75199
$ws.Range("J141").Value = 75199
$ws.Range("L141").Value = 75199
$ws.Range("N141").Value = -85559

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2571
$ws.Range("I7").Value = 1610.8
$ws.Range("J7").Value = 5771.6665
$ws.Range("K7").Value = 1610.8
$ws.Range("L7").Value = 5771.6665
$ws.Range("M7").Value = -1498.8
$ws.Range("N7").Value = -5995.6665
$ws.Range("H40").Value = 50002404
$ws.Range("I40").Value = 55557940
$ws.Range("J40").Value = 2600
$ws.Range("K40").Value = 55557940
$ws.Range("L40").Value = 2600
$ws.Range("M40").Value = -55557804
$ws.Range("N40").Value = -2872
$ws.Range("H61").Value = 2344.6
$ws.Range("I61").Value = 1997.4166
$ws.Range("J61").Value = 3733.3333
$ws.Range("K61").Value = 1997.4166
$ws.Range("L61").Value = 3733.3333
$ws.Range("M61").Value = -1795.4166
$ws.Range("N61").Value = -4137.3333
$ws.Range("H113").Value = 2344.6
$ws.Range("I113").Value = 1997.4166
$ws.Range("J113").Value = 3733.3333
$ws.Range("K113").Value = 1997.4166
$ws.Range("L113").Value = 3733.3333
$ws.Range("M113").Value = 172.5834
$ws.Range("N113").Value = -8073.3333
$ws.Range("H122").Value = 8149655
$ws.Range("I122").Value = 8936694
$ws.Range("J122").Value = 5001500
$ws.Range("K122").Value = 26810082
$ws.Range("L122").Value = 15004500
$ws.Range("M122").Value = -26807632
$ws.Range("N122").Value = -15009400
$ws.Range("H126").Value = 2571
$ws.Range("I126").Value = 1610.8
$ws.Range("J126").Value = 5771.6665
$ws.Range("K126").Value = 4832.4
$ws.Range("L126").Value = 17314.9995
$ws.Range("M126").Value = -2362.4
$ws.Range("N126").Value = -22254.9995
$ws.Range("H141").Value = 79515.625
$ws.Range("J141").Value = 79515.625
$ws.Range("L141").Value = 79515.625
$ws.Range("N141").Value = -89875.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 2505.5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2505.5
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 2505.5
$ws.Range("M29").Value = ""
$ws.Range("N29").Value = -3085.5
$ws.Range("H122").Value = 1470.8
$ws.Range("I122").Value = 1463.5
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 4390.5
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -1940.5
$ws.Range("N122").Value = -9400

